$d = $word.ActiveDocument

# Turn on Track Changes so that each separate insertion step lands in its
# own run (mirroring how the original author's incremental typing produced
# multiple adjacent <w:r> elements). We accept all revisions at the end so
# the final document has no revision marks, just the resulting runs.
$d.TrackRevisions = $true

# --- Paragraph after "How did Mathf.Clamp() help control player boundaries?" ---
$p1 = $d.Paragraphs(33)
$r1 = $p1.Range
$r1.End = $r1.End - 1
$r1.Find.Execute("-", $true, $false, $false, $false, $false, $true, 1, $false, "It constrained the movement between t", 2)

$r1b = $p1.Range
$r1b.End = $r1b.End - 1
$r1b.Collapse(0)
$r1b.InsertAfter("w")

$r1c = $p1.Range
$r1c.End = $r1c.End - 1
$r1c.Collapse(0)
$r1c.InsertAfter("o float values")

$r1d = $p1.Range
$r1d.End = $r1d.End - 1
$r1d.Collapse(0)
$r1d.InsertAfter(".")

# --- Paragraph after "Describe the impact of rotation settings on gameplay?" ---
$p2 = $d.Paragraphs(35)
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Find.Execute("-", $true, $false, $false, $false, $false, $true, 1, $false, "Helped to give the ship a larger area to aim.", 2)

# --- Paragraph after "What adjustments did you make during the tuning process?" ---
$p3 = $d.Paragraphs(37)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Find.Execute("-", $true, $false, $false, $false, $false, $true, 1, $false, "I made the camera further away and more zoomed in as well as pushing the camera up instead of directly behind the ship.", 2)

$d.TrackRevisions = $false
$d.AcceptAllRevisions()
